$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.953.82'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '3.120.50'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.09'
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.31'
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '3.115.88'
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("E10").Value = '  +8.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.76'
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("E12").Value = '  -1.41%  '
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("E14").Value = '  +4.29%  '
$ws.Range("E15").Value = '  -0.83%  '
$ws.Range("D16").Value = '3.636.81'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").Value = '63.847.89'
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.14'
$ws.Range("E18").Value = '  -1.85%  '
$ws.Range("D19").Value = '3.119.70'
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.60'
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.33'
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("E22").Value = '  -0.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("E23").Value = '  -0.94%  '
$ws.Range("E24").Value = '  -3.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.75'
$ws.Range("E25").Value = '  -0.71%  '
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.94'
$ws.Range("E27").Value = '  +7.54%  '
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("E29").Value = '  -1.26%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.00'
$ws.Range("E32").Value = '  -0.48%  '
$ws.Range("E33").Value = '  -2.68%  '
$ws.Range("D34").Value = '0.0₃0867'
$ws.Range("E34").Value = '  +5.38%  '
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.40'
$ws.Range("E37").Value = '  +8.72%  '
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '448.01'
$ws.Range("E40").Value = '  +4.24%  '
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0371'
$ws.Range("E42").Value = '  -0.91%  '
$ws.Range("D43").Value = '2.879.73'
$ws.Range("E43").Value = '  -3.16%  '
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.16'
$ws.Range("E46").Value = '  -1.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.80'
$ws.Range("E47").Value = '  +3.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.57'
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("E51").Value = '  -1.31%  '
